$wb = $excel.ActiveWorkbook

# --- Budget Out sheet: shared string text tweak + amount change ---
$budgetOut = $wb.Worksheets.Item("Budget Out")
$budgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$budgetOut.Range("C9").Value = 93.82

# --- TestRecord sheet: shared string text tweak + date/amount change ---
$testRecord = $wb.Worksheets.Item("TestRecord")
$testRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$testRecord.Range("A10").Value = 43266
$testRecord.Range("B10").Value = 124.74

# --- Expected Out sheet: amount changes (B1 total recalculates automatically) ---
$expectedOut = $wb.Worksheets.Item("Expected Out")
$expectedOut.Range("B9").Value = 1351.76
$expectedOut.Range("B11").Value = 431.62
